$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-03-02 14:54:17"
$wsZhCn.Range("G4").Value = "2016-03-02 14:55:01"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-03-02 14:54:27"
$wsDeDe.Range("G4").Value = "2016-03-02 14:55:32"
